$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do- FY16 Release")

# --- Update a few existing cells ---
$ws.Range("A8").Value = "Future"
$ws.Range("A13").Value = "Testing"
$ws.Range("B18").Value = "Linear shading option for PV self-shading model (and upgrade project files!)"
$ws.Range("C20").Value = "Janine/Paul"

# --- Insert 16 new rows before row 22 for new TODO items ---
$ws.Rows("22:37").Insert()

$ws.Range("A22").Value = "Not Done"
$ws.Range("B22").Value = "Battery automatic dispatch and testing, project file upgrader, etc"
$ws.Range("C22").Value = "Nick"

$ws.Range("A23").Value = "Testing"
$ws.Range("B23").Value = "Finish testing of POA model, project upgrade, documentation"
$ws.Range("C23").Value = "Janine"

$ws.Range("A24").Value = "Not Done"
$ws.Range("B24").Value = "Updates to utility rate model metering options and input/output structures"
$ws.Range("C24").Value = "Steve/Paul"

$ws.Range("A25").Value = "Not Done"
$ws.Range("B25").Value = "Proxy autodetection on Windows "
$ws.Range("C25").Value = "Aron"

$ws.Range("A26").Value = "Testing"
$ws.Range("B26").Value = "Validation and fixes to 3D shading database"
$ws.Range("C26").Value = "Steve/Sara"

$ws.Range("A27").Value = "Not Done"
$ws.Range("B27").Value = "Remove webkitgtk dependencies on Linux OS, rebuild on CentOS 6.4 as standard"
$ws.Range("C27").Value = "Aron"

$ws.Range("A28").Value = "Done"
$ws.Range("B28").Value = "Wind wizard"
$ws.Range("C28").Value = "Janine"

$ws.Range("A29").Value = "Testing"
$ws.Range("B29").Value = "User defined power cycle input option"
$ws.Range("C29").Value = "Ty"

$ws.Range("A30").Value = "Testing"
$ws.Range("B30").Value = "Molten salt power tower model with new controller framework"
$ws.Range("C30").Value = "Ty"

$ws.Range("A31").Value = "Testing"
$ws.Range("B31").Value = "Dispatch optimization for power tower model"
$ws.Range("C31").Value = "Mike"

$ws.Range("A32").Value = "Not Done"
$ws.Range("B32").Value = "Merge updates from SolarPilot standalone version into SAM, update defaults"
$ws.Range("C32").Value = "Mike"

$ws.Range("A33").Value = "Testing"
$ws.Range("B33").Value = "Check s3d file format reader: does it read old .s3d files OK with group property removed/renamed? Notice in project file upgrader"
$ws.Range("C33").Value = "Steve"

$ws.Range("A34").Value = "Not Done"
$ws.Range("B34").Value = "LK 'global' variable space keyword"
$ws.Range("C34").Value = "Aron"

$ws.Range("A35").Value = "Not Done"
$ws.Range("B35").Value = "Update LK documentation in line with language changes in VM"
$ws.Range("C35").Value = "Aron"

$ws.Range("A36").Value = "Future"
$ws.Range("B36").Value = "Checkbox on tower page to update inputs from results based on a script.  Add a post-simulation script framework"
$ws.Range("C36").Value = "Aron/Ty"

$ws.Range("A37").Value = "Future"
$ws.Range("B37").Value = "Remove HCPV model and augment simple efficiency model accordingly"
$ws.Range("C37").Value = "Aron"

# --- Update selection to match the target view state ---
$ws.Range("C13").Select()
